# Update player data, activity logs, and war records; add new input files and adjust strike details
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update player names in column B (rows 2-31) ---
$ws.Cells.Item(2, 2).Value = "Satan"
$ws.Cells.Item(3, 2).Value = "Kaselcap"
$ws.Cells.Item(4, 2).Value = "BumblinMumbler"
$ws.Cells.Item(5, 2).Value = "Sned"
$ws.Cells.Item(6, 2).Value = "Smitty™"
$ws.Cells.Item(7, 2).Value = "Protips"
$ws.Cells.Item(8, 2).Value = "katsu"
$ws.Cells.Item(9, 2).Value = "pg"
$ws.Cells.Item(10, 2).Value = "Big Daddy T"
$ws.Cells.Item(11, 2).Value = "Vojt"
$ws.Cells.Item(12, 2).Value = "K.L.A.U.S"
$ws.Cells.Item(13, 2).Value = "Anas"
$ws.Cells.Item(14, 2).Value = "Az7777"
$ws.Cells.Item(15, 2).Value = "Hadez"
$ws.Cells.Item(16, 2).Value = "Ascended"
$ws.Cells.Item(17, 2).Value = "Baleus"
$ws.Cells.Item(18, 2).Value = "YouAreMyBreh"
$ws.Cells.Item(19, 2).Value = "Rod"
$ws.Cells.Item(20, 2).Value = "xHead_Bangerx"
$ws.Cells.Item(21, 2).Value = "Mythos"
$ws.Cells.Item(22, 2).Value = "Plantos"
$ws.Cells.Item(23, 2).Value = "ImagineWaggons"
$ws.Cells.Item(24, 2).Value = "Luke"
$ws.Cells.Item(25, 2).Value = "shadow"
$ws.Cells.Item(26, 2).Value = "Welli"
$ws.Cells.Item(27, 2).Value = "Motz"
$ws.Cells.Item(28, 2).Value = "LOGAN911"
$ws.Cells.Item(29, 2).Value = "potatoes"
$ws.Cells.Item(30, 2).Value = "DNG"
$ws.Cells.Item(31, 2).Value = "Asrar"

# --- Clear all existing war-strike cells (D:J) before rebuilding them to match new records ---
$ws.Range("D2:J31").Clear()

# --- Style donor cells that persist in column B throughout (never cleared): ---
# style index "2" (rows 2-13 group) lives on B2; style index "4" (rows 14-31 group) lives on B14
$style2Donor = $ws.Range("B2")
$style4Donor = $ws.Range("B14")

# --- Rebuild war-strike ("fwa") marks and blank placeholder cells per row ---
# Row 2: Satan
$style4Donor.Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = "fwa"
$style4Donor.Copy()
$ws.Range("E2").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("F2").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("G2").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("H2").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("I2").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("J2").PasteSpecial(-4122)

# Row 3: Kaselcap
$style2Donor.Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "fwa"
$style2Donor.Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "fwa"

# Row 10: Big Daddy T
$style2Donor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = "fwa"

# Row 12: K.L.A.U.S
$style4Donor.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("E12").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("F12").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("G12").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("H12").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("I12").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("J12").PasteSpecial(-4122)

# Row 15: Hadez
$style4Donor.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("F15").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("G15").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("I15").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("J15").PasteSpecial(-4122)

# Row 18: YouAreMyBreh
$style2Donor.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "fwa"
$style2Donor.Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "fwa"

# Row 21: Mythos
$style4Donor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Value = "fwa"
$style4Donor.Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "fwa"
$style4Donor.Copy()
$ws.Range("F21").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("G21").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("H21").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("I21").PasteSpecial(-4122)
$style4Donor.Copy()
$ws.Range("J21").PasteSpecial(-4122)

# Row 28: LOGAN911
$style2Donor.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = "fwa"
$style2Donor.Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "fwa"
$style2Donor.Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = "fwa"
$style2Donor.Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = "fwa"
$style2Donor.Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = "fwa"

$ws.Application.CutCopyMode = 0
